$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1135
$ws1.Range("F5").Value = 185
$ws1.Range("F6").Value = 5
$ws1.Range("F8").Value = 253
$ws1.Range("F11").Value = 16
$ws1.Range("F15").Value = 12912
$ws1.Range("F16").Value = 9
$ws1.Range("F17").Value = 5302

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 143

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1135
$ws4.Range("F5").Value = 185
$ws4.Range("F6").Value = 5
$ws4.Range("F8").Value = 253
$ws4.Range("F11").Value = 16
$ws4.Range("F15").Value = 12912
$ws4.Range("F16").Value = 143
$ws4.Range("F18").Value = 9
$ws4.Range("F19").Value = 5302
